$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 82, shifting the existing rows
# 82-193 down to 84-195 (this also extends the sheet dimension to R195).
$ws.Rows("82:83").Insert()

# Fill in the new row 82 with its data.
$ws.Cells.Item(82, 1).Value2 = 10
$ws.Cells.Item(82, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(82, 3).Value2 = "La Araucanía"
$ws.Cells.Item(82, 4).Value2 = 44467
$ws.Cells.Item(82, 5).Value2 = 9
$ws.Cells.Item(82, 6).Value2 = 100112044
$ws.Cells.Item(82, 7).Value2 = "Perejil"
$ws.Cells.Item(82, 8).Value2 = "Sin especificar"
$ws.Cells.Item(82, 9).Value2 = "Primera"
$ws.Cells.Item(82, 10).Value2 = 60
$ws.Cells.Item(82, 11).Value2 = 3000
$ws.Cells.Item(82, 12).Value2 = 4000
$ws.Cells.Item(82, 13).Value2 = 3500
$ws.Cells.Item(82, 14).Value2 = "$/docena de atados (3 kilos)"
$ws.Cells.Item(82, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(82, 16).Value2 = 1167
$ws.Cells.Item(82, 17).Value2 = 3
$ws.Cells.Item(82, 18).Value2 = "Hortaliza"

# Fill in the new row 83 with its data.
$ws.Cells.Item(83, 1).Value2 = 10
$ws.Cells.Item(83, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(83, 3).Value2 = "La Araucanía"
$ws.Cells.Item(83, 4).Value2 = 44467
$ws.Cells.Item(83, 5).Value2 = 9
$ws.Cells.Item(83, 6).Value2 = 100112044
$ws.Cells.Item(83, 7).Value2 = "Perejil"
$ws.Cells.Item(83, 8).Value2 = "Sin especificar"
$ws.Cells.Item(83, 9).Value2 = "Primera"
$ws.Cells.Item(83, 10).Value2 = 30
$ws.Cells.Item(83, 11).Value2 = 3300
$ws.Cells.Item(83, 12).Value2 = 3300
$ws.Cells.Item(83, 13).Value2 = 3300
$ws.Cells.Item(83, 14).Value2 = "$/docena de atados (3 kilos)"
$ws.Cells.Item(83, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(83, 16).Value2 = 1100
$ws.Cells.Item(83, 17).Value2 = 3
$ws.Cells.Item(83, 18).Value2 = "Hortaliza"

# Make sure the date cells keep the same number format used throughout
# column D (inherited automatically from the row-insert, but set explicitly
# to be safe).
$ws.Range("D82:D83").NumberFormat = $ws.Range("D84").NumberFormat
